# Fruta / hortaliza, semanal
# A new weekly price record is inserted as the new row 43 on the
# "Frambuesa" sheet, pushing the previously existing rows 43-61 down to
# rows 44-62 (dimension grows from A1:T61 to A1:T62).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 43; Excel shifts rows
# 43:61 down to 44:62 and copies the row-above formatting (keeps the
# date style on column D) automatically.
$ws.Rows.Item(43).Insert()

# Populate the newly inserted row 43 with the new weekly record.
$ws.Cells.Item(43, 1).Value = 9
$ws.Cells.Item(43, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(43, 3).Value = 'Metropolitana'
$ws.Cells.Item(43, 4).Value = 44559
$ws.Cells.Item(43, 5).Value = 13
$ws.Cells.Item(43, 6).Value = 'Fruta'
$ws.Cells.Item(43, 7).Value = 100101
$ws.Cells.Item(43, 8).Value = 'Berries'
$ws.Cells.Item(43, 9).Value = 100101004
$ws.Cells.Item(43, 10).Value = 'Frambuesa'
$ws.Cells.Item(43, 11).Value = 'Sin especificar'
$ws.Cells.Item(43, 12).Value = 'Primera'
$ws.Cells.Item(43, 13).Value = 450
$ws.Cells.Item(43, 14).Value = 8000
$ws.Cells.Item(43, 15).Value = 8000
$ws.Cells.Item(43, 16).Value = 8000
$ws.Cells.Item(43, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(43, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(43, 19).Value = 4000
$ws.Cells.Item(43, 20).Value = 2
